$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.181.91"
$ws.Range("E2").Value = "  +0.64%  "

$ws.Range("D3").Value = "2.415.85"
$ws.Range("E3").Value = "  +1.37%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'560.84"
$ws.Range("E5").Value = "  +1.42%  "

$ws.Range("D6").Value = "'142.83"
$ws.Range("E6").Value = "  +0.93%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'0.531"
$ws.Range("E8").Value = "  +1.02%  "

$ws.Range("D9").Value = "2.411.12"
$ws.Range("E9").Value = "  +0.92%  "

$ws.Range("E10").Value = "  +0.23%  "

$ws.Range("E11").Value = "  -2.02%  "

$ws.Range("D12").Value = "'5.32"
$ws.Range("E12").Value = "  -0.99%  "

$ws.Range("D13").Value = "'0.352"
$ws.Range("E13").Value = "  -0.68%  "

$ws.Range("D14").Value = "'25.67"
$ws.Range("E14").Value = "  -0.46%  "

$ws.Range("D15").Value = "'0.0000175"
$ws.Range("E15").Value = "  -0.30%  "

$ws.Range("D16").Value = "2.859.88"
$ws.Range("E16").Value = "  +1.55%  "

$ws.Range("D17").Value = "62.070.47"
$ws.Range("E17").Value = "  +0.56%  "

$ws.Range("D18").Value = "2.417.55"
$ws.Range("E18").Value = "  +1.30%  "

$ws.Range("D19").Value = "'11.26"
$ws.Range("E19").Value = "  +2.16%  "

$ws.Range("D20").Value = "'4.18"
$ws.Range("E20").Value = "  -0.26%  "

$ws.Range("D21").Value = "'322.87"
$ws.Range("E21").Value = "  -0.34%  "

$ws.Range("D22").Value = "'6.80"
$ws.Range("E22").Value = "  +1.65%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").Value = "'65.77"
$ws.Range("E24").Value = "  +2.08%  "

$ws.Range("D25").Value = "'1.71"
$ws.Range("E25").Value = "  -4.32%  "

$ws.Range("E26").Value = "  +0.34%  "

$ws.Range("D27").Value = "'578.21"
$ws.Range("E27").Value = "  +6.61%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.536.92"
$ws.Range("E28").Value = "  +1.37%  "

$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.45%  "

$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0941"
$ws.Range("E30").Value = "  +1.65%  "

$ws.Range("D31").Value = "'8.21"
$ws.Range("E31").Value = "  -1.15%  "

$ws.Range("D32").Value = "'1.44"
$ws.Range("E32").Value = "  +0.85%  "

$ws.Range("D33").Value = "'0.149"
$ws.Range("E33").Value = "  +0.28%  "

$ws.Range("D34").Value = "'1.87"
$ws.Range("E34").Value = "  +0.82%  "

$ws.Range("D35").Value = "'1.54"
$ws.Range("E35").Value = "  -0.27%  "

$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("D37").Value = "'5.63"
$ws.Range("E37").Value = "  -2.14%  "

$ws.Range("D38").Value = "'4.73"
$ws.Range("E38").Value = "  -0.93%  "

$ws.Range("D39").Value = "'0.383"
$ws.Range("E39").Value = "  +0.38%  "

$ws.Range("D40").Value = "'152.25"
$ws.Range("E40").Value = "  +3.60%  "

$ws.Range("D41").Value = "'18.61"
$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").Value = "'1.81"
$ws.Range("E42").Value = "  -6.74%  "

$ws.Range("E43").Value = "  -0.29%  "

$ws.Range("D44").Value = "'2.31"
$ws.Range("E44").Value = "  +2.47%  "

$ws.Range("D45").Value = "'149.04"
$ws.Range("E45").Value = "  +0.24%  "

$ws.Range("D46").Value = "'3.65"
$ws.Range("E46").Value = "  +0.74%  "

$ws.Range("D47").Value = "'0.0535"
$ws.Range("E47").Value = "  +0.84%  "

$ws.Range("D48").Value = "'20.07"
$ws.Range("E48").Value = "  -0.66%  "

$ws.Range("D49").Value = "'0.593"
$ws.Range("E49").Value = "  +1.17%  "

$ws.Range("D50").Value = "'0.0917"
$ws.Range("E50").Value = "  +1.02%  "

$ws.Range("D51").Value = "'0.0228"
$ws.Range("E51").Value = "  +1.07%  "
